# The "7-11-21" $1.00 game (row 9) was removed from the NY scraper sheet.
# Deleting the entire row shifts every row below it up by one and Excel
# automatically drops the now-unused "7-11-21" shared string / shrinks the
# used range from A1:F78 to A1:F77.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9:F9").EntireRow.Delete()
